# Update sample excel file for the excel-test-results-plugin:
# add a new "Error" column (E) with an error message for the failing
# "Add two numbers (outline), negative" test row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "Error"

# New error detail for row 6 (the Failed test row)
$ws.Range("E6").Value = "An OutOfRange error was thrown"

# Move/update the active selection to match the edited workbook (E7)
$ws.Range("E7").Select() | Out-Null
